# Apply updated cryptocurrency price/volume data (and the Kaspa/Cosmos row swap)
# to match the refreshed source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.955.76'
$ws.Range('E2').Value = '  +0.40%  '

$ws.Range('D3').Value = '3.433.42'
$ws.Range('E3').Value = '  +1.27%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '''409.06'
$ws.Range('E5').Value = '  +0.93%  '

$ws.Range('D6').Value = '''128.78'
$ws.Range('E6').Value = '  -3.66%  '

$ws.Range('D7').Value = '''0.623'
$ws.Range('E7').Value = '  +5.61%  '

$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('D9').Value = '''0.742'
$ws.Range('E9').Value = '  +10.75%  '

$ws.Range('D10').Value = '''0.143'
$ws.Range('E10').Value = '  +17.85%  '

$ws.Range('D11').Value = '''42.78'
$ws.Range('E11').Value = '  +0.71%  '

$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('D13').Value = '3.970.97'
$ws.Range('E13').Value = '  +1.45%  '

$ws.Range('D14').Value = '''21.34'
$ws.Range('E14').Value = '  +7.91%  '

$ws.Range('D15').Value = '''8.92'
$ws.Range('E15').Value = '  +5.84%  '

$ws.Range('D16').Value = '''0.0000209'
$ws.Range('E16').Value = '  +62.77%  '

$ws.Range('D17').Value = '3.442.72'
$ws.Range('E17').Value = '  +1.57%  '

$ws.Range('D18').Value = '''12.41'
$ws.Range('E18').Value = '  +12.52%  '

$ws.Range('E19').Value = '  +5.40%  '

$ws.Range('D20').Value = '61.910.18'
$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('D21').Value = '''401.90'
$ws.Range('E21').Value = '  +27.64%  '

$ws.Range('D22').Value = '''90.01'
$ws.Range('E22').Value = '  +5.26%  '

$ws.Range('D23').Value = '''3.19'
$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').Value = '''13.42'
$ws.Range('E24').Value = '  +5.32%  '

$ws.Range('E25').Value = '  +3.19%  '

$ws.Range('D26').Value = '''33.21'
$ws.Range('E26').Value = '  +12.38%  '

$ws.Range('D27').Value = '''8.71'
$ws.Range('E27').Value = '  +4.74%  '

$ws.Range('D28').Value = '''4.79'
$ws.Range('E28').Value = '  +0.10%  '

$ws.Range('E29').Value = '  -1.00%  '

$ws.Range('D30').Value = '''2.73'
$ws.Range('E30').Value = '  +2.72%  '

$ws.Range('E31').Value = '  +2.73%  '

$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '''0.172'
$ws.Range('E32').Value = '  +0.47%  '

$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''11.90'
$ws.Range('E33').Value = '  +4.99%  '

$ws.Range('D34').Value = '''43.44'
$ws.Range('E34').Value = '  +4.83%  '

$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('D36').Value = '''0.0504'
$ws.Range('E36').Value = '  +4.99%  '

$ws.Range('D37').Value = '''54.28'
$ws.Range('E37').Value = '  +4.77%  '

$ws.Range('E38').Value = '  +0.17%  '

$ws.Range('D39').Value = '''3.39'
$ws.Range('E39').Value = '  -0.89%  '

$ws.Range('E40').Value = '  -0.74%  '

$ws.Range('E41').Value = '  +6.30%  '

$ws.Range('E42').Value = '  +6.34%  '

$ws.Range('D43').Value = '''141.61'
$ws.Range('E43').Value = '  +1.69%  '

$ws.Range('D44').Value = '''1.98'
$ws.Range('E44').Value = '  +0.06%  '

$ws.Range('D45').Value = '''4.04'
$ws.Range('E45').Value = '  +1.85%  '

$ws.Range('D46').Value = '''2.40'
$ws.Range('E46').Value = '  +7.93%  '

$ws.Range('D47').Value = '''16.73'
$ws.Range('E47').Value = '  +0.62%  '

$ws.Range('D48').Value = '''21.80'
$ws.Range('E48').Value = '  +2.43%  '

$ws.Range('D49').Value = '2.121.20'
$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('E50').Value = '  +14.90%  '

$ws.Range('D51').Value = '''0.0376'
$ws.Range('E51').Value = '  +7.61%  '
